$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.037623761840726
$ws.Range("D2").Value = 1.043173887947079
$ws.Range("E2").Value = 1.045526345577224
$ws.Range("F2").Value = 1.054336280728075
$ws.Range("I2").Value = 1.034503794856443
$ws.Range("J2").Value = 1.042725623761427
$ws.Range("K2").Value = 1.045948498789237
$ws.Range("L2").Value = 1.048294338033103
$ws.Range("M2").Value = 1.057079771785203

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038577805455057
$ws.Range("D3").Value = 1.043883638273902
$ws.Range("E3").Value = 1.046358070655741
$ws.Range("F3").Value = 1.055213040439652
$ws.Range("I3").Value = 1.034633692900023
$ws.Range("J3").Value = 1.043324100484633
$ws.Range("K3").Value = 1.046469495277279
$ws.Range("L3").Value = 1.048937468656444
$ws.Range("M3").Value = 1.057769592038708

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039195782692573
$ws.Range("D4").Value = 1.044343245103835
$ws.Range("E4").Value = 1.046897167687353
$ws.Range("F4").Value = 1.055781259504697
$ws.Range("I4").Value = 1.034716500506094
$ws.Range("J4").Value = 1.043711384291605
$ws.Range("K4").Value = 1.046806292092814
$ws.Range("L4").Value = 1.049353883366778
$ws.Range("M4").Value = 1.058216207624401

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039455734109401
$ws.Range("D5").Value = 1.044536546244919
$ws.Range("E5").Value = 1.047124021658704
$ws.Range("F5").Value = 1.05602035189166
$ws.Range("I5").Value = 1.03475101415101
$ws.Range("J5").Value = 1.043874204300288
$ws.Range("K5").Value = 1.046947803017557
$ws.Range("L5").Value = 1.049529006600434
$ws.Range("M5").Value = 1.058404024676444

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039499390080592
$ws.Range("D6").Value = 1.04456900714488
$ws.Range("E6").Value = 1.047162124152799
$ws.Range("F6").Value = 1.056060509003967
$ws.Range("I6").Value = 1.034756791612421
$ws.Range("J6").Value = 1.043901542809124
$ws.Range("K6").Value = 1.046971558689196
$ws.Range("L6").Value = 1.04955841419039
$ws.Range("M6").Value = 1.058435563456354

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039199255573299
$ws.Range("D7").Value = 1.044345827681237
$ws.Range("E7").Value = 1.046900198068221
$ws.Range("F7").Value = 1.055784453433614
$ws.Range("I7").Value = 1.034716962853423
$ws.Range("J7").Value = 1.043713559878779
$ws.Range("K7").Value = 1.046808183278491
$ws.Range("L7").Value = 1.04935622312801
$ws.Range("M7").Value = 1.058218717011554

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037946050635842
$ws.Range("D8").Value = 1.043413677786369
$ws.Range("E8").Value = 1.045807240517035
$ws.Range("F8").Value = 1.054632399088222
$ws.Range("I8").Value = 1.034547951962691
$ws.Range("J8").Value = 1.042927875015869
$ws.Range("K8").Value = 1.046124638104078
$ws.Range("L8").Value = 1.048511631177376
$ws.Range("M8").Value = 1.057312845991052

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035742751723985
$ws.Range("D9").Value = 1.041773867822332
$ws.Range("E9").Value = 1.04388838787467
$ws.Range("F9").Value = 1.052609273327367
$ws.Range("I9").Value = 1.034240623171541
$ws.Range("J9").Value = 1.041543672589867
$ws.Range("K9").Value = 1.044917724830253
$ws.Range("L9").Value = 1.047025445638035
$ws.Range("M9").Value = 1.055718604747989

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034277321333194
$ws.Range("D10").Value = 1.040682608668339
$ws.Range("E10").Value = 1.042614003290575
$ws.Range("F10").Value = 1.051265286297231
$ws.Range("I10").Value = 1.034029378329756
$ws.Range("J10").Value = 1.040621122499976
$ws.Range("K10").Value = 1.044111559838688
$ws.Range("L10").Value = 1.046036136064702
$ws.Range("M10").Value = 1.054657215980735

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033643602732423
$ws.Range("D11").Value = 1.040210562984776
$ws.Range("E11").Value = 1.042063350354265
$ws.Range("F11").Value = 1.050684475153296
$ws.Range("I11").Value = 1.033936406032986
$ws.Range("J11").Value = 1.040221721730945
$ws.Range("K11").Value = 1.043762127147613
$ws.Range("L11").Value = 1.045608120649844
$ws.Range("M11").Value = 1.054197981147683

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033408336081603
$ws.Range("D12").Value = 1.040035297397678
$ws.Range("E12").Value = 1.041858989750477
$ws.Range("F12").Value = 1.050468909544062
$ws.Range("I12").Value = 1.03390164690701
$ws.Range("J12").Value = 1.040073378046883
$ws.Range("K12").Value = 1.043632279680375
$ws.Range("L12").Value = 1.045449192296452
$ws.Range("M12").Value = 1.054027455344625

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033458795918378
$ws.Range("D13").Value = 1.040072889123457
$ws.Range("E13").Value = 1.041902817785897
$ws.Range("F13").Value = 1.050515141221266
$ws.Range("I13").Value = 1.033909113033559
$ws.Range("J13").Value = 1.040105197730541
$ws.Range("K13").Value = 1.043660134764714
$ws.Range("L13").Value = 1.045483280430841
$ws.Range("M13").Value = 1.054064031222057

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033624152968432
$ws.Range("D14").Value = 1.040196073974365
$ws.Range("E14").Value = 1.042046454235312
$ws.Range("F14").Value = 1.05066665288386
$ws.Range("I14").Value = 1.03393353742098
$ws.Range("J14").Value = 1.040209459353854
$ws.Range("K14").Value = 1.043751394981884
$ws.Range("L14").Value = 1.045594982440707
$ws.Range("M14").Value = 1.054183884312797

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033726051455081
$ws.Range("D15").Value = 1.040271981973077
$ws.Range("E15").Value = 1.042134976815017
$ws.Range("F15").Value = 1.050760027259978
$ws.Range("I15").Value = 1.033948556284461
$ws.Range("J15").Value = 1.04027369994384
$ws.Range("K15").Value = 1.043807616477605
$ws.Range("L15").Value = 1.045663813143477
$ws.Range("M15").Value = 1.05425773701473

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034319396570442
$ws.Range("D16").Value = 1.040713946964288
$ws.Range("E16").Value = 1.042650572967727
$ws.Range("F16").Value = 1.051303857063402
$ws.Range("I16").Value = 1.03403551700017
$ws.Range("J16").Value = 1.040647630972119
$ws.Range("K16").Value = 1.044134743094056
$ws.Range("L16").Value = 1.046064549772443
$ws.Range("M16").Value = 1.054687701447442

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034691807242931
$ws.Range("D17").Value = 1.040991308672883
$ws.Range("E17").Value = 1.04297430553285
$ws.Range("F17").Value = 1.051645294585855
$ws.Range("I17").Value = 1.034089663461308
$ws.Range("J17").Value = 1.04088220741392
$ws.Range("K17").Value = 1.0443398459742
$ws.Range("L17").Value = 1.046316019092896
$ws.Range("M17").Value = 1.054957502374002

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034909107366347
$ws.Range("D18").Value = 1.041153134957077
$ws.Range("E18").Value = 1.043163245466836
$ws.Range("F18").Value = 1.051844559597254
$ws.Range("I18").Value = 1.03412110122773
$ws.Range("J18").Value = 1.041019038516217
$ws.Range("K18").Value = 1.044459444400752
$ws.Range("L18").Value = 1.046462731754294
$ws.Range("M18").Value = 1.055114906707487

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034983214506783
$ws.Range("D19").Value = 1.041208321274221
$ws.Range("E19").Value = 1.043227688086984
$ws.Range("F19").Value = 1.051912522472976
$ws.Range("I19").Value = 1.034131796090708
$ws.Range("J19").Value = 1.041065695477547
$ws.Range("K19").Value = 1.044500218466963
$ws.Range("L19").Value = 1.046512762856448
$ws.Range("M19").Value = 1.055168583246416

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034651842894129
$ws.Range("D20").Value = 1.040961545628258
$ws.Range("E20").Value = 1.042939560475299
$ws.Range("F20").Value = 1.051608650152411
$ws.Range("I20").Value = 1.034083869047056
$ws.Range("J20").Value = 1.040857038887721
$ws.Range("K20").Value = 1.044317843948016
$ws.Range("L20").Value = 1.046289035200013
$ws.Range("M20").Value = 1.054928551759577

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033575455988744
$ws.Range("D21").Value = 1.040159797072825
$ws.Range("E21").Value = 1.042004151998753
$ws.Range("F21").Value = 1.050622031672804
$ws.Range("I21").Value = 1.033926351256276
$ws.Range("J21").Value = 1.040178756576888
$ws.Range("K21").Value = 1.043724522565368
$ws.Range("L21").Value = 1.045562087437474
$ws.Range("M21").Value = 1.054148589050581

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032899410573527
$ws.Range("D22").Value = 1.039656130211825
$ws.Range("E22").Value = 1.041417045004608
$ws.Range("F22").Value = 1.050002710155277
$ws.Range("I22").Value = 1.033826011633792
$ws.Range("J22").Value = 1.039752360908698
$ws.Range("K22").Value = 1.043351174071988
$ws.Range("L22").Value = 1.045105348983195
$ws.Range("M22").Value = 1.053658511304659

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033257726258911
$ws.Range("D23").Value = 1.039923092802685
$ws.Range("E23").Value = 1.041728184196542
$ws.Range("F23").Value = 1.050330928436362
$ws.Range("I23").Value = 1.033879326785279
$ws.Range("J23").Value = 1.039978394568546
$ws.Range("K23").Value = 1.043549121578246
$ws.Range("L23").Value = 1.045347443707287
$ws.Range("M23").Value = 1.05391828040816

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034669900818232
$ws.Range("D24").Value = 1.040974994125366
$ws.Range("E24").Value = 1.042955259924696
$ws.Range("F24").Value = 1.051625207853971
$ws.Range("I24").Value = 1.034086487741363
$ws.Range("J24").Value = 1.040868411440959
$ws.Range("K24").Value = 1.044327785823653
$ws.Range("L24").Value = 1.046301227952111
$ws.Range("M24").Value = 1.054941633191105

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03631175678202
$ws.Range("D25").Value = 1.042197461599447
$ws.Range("E25").Value = 1.044383609043188
$ws.Range("F25").Value = 1.053131467310276
$ws.Range("I25").Value = 1.034321198412053
$ws.Range("J25").Value = 1.041901482660622
$ws.Range("K25").Value = 1.045230020283457
$ws.Range("L25").Value = 1.047409404594498
$ws.Range("M25").Value = 1.056130506677703
